# CLEVER LULUCF assumptions for BAU scenario + small changes
# -------------------------------------------------------------
# Fills in the "Input_File" (column B) entries on the COUNTRIES sheet for
# every country that was still missing one, following the "inputsXX"
# naming convention already used for BE / DE / FR / NL / GB / EU. Greece
# (row 11, ISO "EL") gets "inputsGR" (its input file keys off the ISO
# "GR" prefix rather than the Eurostat "EL" code). Cyprus (row 6) and
# Malta (row 22) intentionally keep no Input_File value, matching the
# author's commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("COUNTRIES")
$ws.Activate()

$inputFiles = @(
    @{ Row = 2;  Value = "inputsAT" },   # AT - Austria
    @{ Row = 4;  Value = "inputsBG" },   # BG - Bulgaria
    @{ Row = 5;  Value = "inputsCH" },   # CH - Switzerland
    @{ Row = 7;  Value = "inputsCZ" },   # CZ - Czech Republic
    @{ Row = 9;  Value = "inputsDK" },   # DK - Denmark
    @{ Row = 10; Value = "inputsEE" },   # EE - Estonia
    @{ Row = 11; Value = "inputsGR" },   # EL - Greece
    @{ Row = 12; Value = "inputsES" },   # ES - Spain
    @{ Row = 13; Value = "inputsFI" },   # FI - Finland
    @{ Row = 15; Value = "inputsHR" },   # HR - Croatia
    @{ Row = 16; Value = "inputsHU" },   # HU - Hungary
    @{ Row = 17; Value = "inputsIE" },   # IE - Ireland
    @{ Row = 18; Value = "inputsIT" },   # IT - Italy
    @{ Row = 19; Value = "inputsLT" },   # LT - Lithuania
    @{ Row = 20; Value = "inputsLU" },   # LU - Luxembourg
    @{ Row = 21; Value = "inputsLV" },   # LV - Latvia
    @{ Row = 24; Value = "inputsNO" },   # NO - Norway
    @{ Row = 25; Value = "inputsPL" },   # PL - Poland
    @{ Row = 26; Value = "inputsPT" },   # PT - Portugal
    @{ Row = 27; Value = "inputsRO" },   # RO - Romania
    @{ Row = 28; Value = "inputsSE" },   # SE - Sweden
    @{ Row = 29; Value = "inputsSI" },   # SI - Slovenia
    @{ Row = 30; Value = "inputsSK" }    # SK - Slovakia
)

foreach ($entry in $inputFiles) {
    $ws.Cells.Item($entry.Row, 2).Value = $entry.Value
}

# Move the selection on the COUNTRIES sheet from C32 to B2 (matches the
# updated sheetView active cell in the saved workbook).
[void]$ws.Range("B2").Select()
